# Update res_line/pl_mw.xlsx numeric results (rows 2-25, columns B,C,D,E,F,H,I,K,L)
# for the "case with 380 kV" re-run.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 2).Value = 0.78016879002368
$ws.Cells.Item(2, 3).Value = 0.0416085903975727
$ws.Cells.Item(2, 4).Value = 0.1433421750881294
$ws.Cells.Item(2, 5).Value = 0.06113185119997766
$ws.Cells.Item(2, 6).Value = 2.630708748447375
$ws.Cells.Item(2, 8).Value = 0.07973214163530429
$ws.Cells.Item(2, 9).Value = 2.147253301571425
$ws.Cells.Item(2, 11).Value = 0.6164355302456386
$ws.Cells.Item(2, 12).Value = 0.2456817495573063

$ws.Cells.Item(3, 2).Value = 0.754189425547878
$ws.Cells.Item(3, 3).Value = 0.03620255580442233
$ws.Cells.Item(3, 4).Value = 0.1426570395147166
$ws.Cells.Item(3, 5).Value = 0.0606836595636544
$ws.Cells.Item(3, 6).Value = 2.577096267988665
$ws.Cells.Item(3, 8).Value = 0.07973214163530429
$ws.Cells.Item(3, 9).Value = 2.114043808880012
$ws.Cells.Item(3, 11).Value = 0.5853739611381172
$ws.Cells.Item(3, 12).Value = 0.238526883287733

$ws.Cells.Item(4, 2).Value = 0.7388538796103887
$ws.Cells.Item(4, 3).Value = 0.03288011789024381
$ws.Cells.Item(4, 4).Value = 0.1422241723565563
$ws.Cells.Item(4, 5).Value = 0.06043441177345876
$ws.Cells.Item(4, 6).Value = 2.545016477067875
$ws.Cells.Item(4, 8).Value = 0.07973214163530429
$ws.Cells.Item(4, 9).Value = 2.094205655803293
$ws.Cells.Item(4, 11).Value = 0.5667910953249873
$ws.Cells.Item(4, 12).Value = 0.2342796997542109

$ws.Cells.Item(5, 2).Value = 0.7327593670989927
$ws.Cells.Item(5, 3).Value = 0.03152531484499832
$ws.Cells.Item(5, 4).Value = 0.1420446831963744
$ws.Cells.Item(5, 5).Value = 0.06033936952589869
$ws.Cells.Item(5, 6).Value = 2.532153828972497
$ws.Cells.Item(5, 8).Value = 0.07973214163530429
$ws.Cells.Item(5, 9).Value = 2.086259911447911
$ws.Cells.Item(5, 11).Value = 0.5593412042660191
$ws.Cells.Item(5, 12).Value = 0.2325855804218264

$ws.Cells.Item(6, 2).Value = 0.7317567298548511
$ws.Cells.Item(6, 3).Value = 0.03130029460452022
$ws.Cells.Item(6, 4).Value = 0.1420146916063523
$ws.Cells.Item(6, 5).Value = 0.06032398226010471
$ws.Cells.Item(6, 6).Value = 2.530030666377584
$ws.Cells.Item(6, 8).Value = 0.07973214163530429
$ws.Cells.Item(6, 9).Value = 2.084948875163704
$ws.Cells.Item(6, 11).Value = 0.5581115653199618
$ws.Cells.Item(6, 12).Value = 0.2323064847839049

$ws.Cells.Item(7, 2).Value = 0.7387710599409161
$ws.Cells.Item(7, 3).Value = 0.03286185024863642
$ws.Cells.Item(7, 4).Value = 0.1422217642593075
$ws.Cells.Item(7, 5).Value = 0.06043310356152176
$ws.Cells.Item(7, 6).Value = 2.544842157063783
$ws.Cells.Item(7, 8).Value = 0.07973214163530429
$ws.Cells.Item(7, 9).Value = 2.094097936746337
$ws.Cells.Item(7, 11).Value = 0.5666901266277193
$ws.Cells.Item(7, 12).Value = 0.2342567039734718

$ws.Cells.Item(8, 2).Value = 0.7710832498189859
$ws.Cells.Item(8, 3).Value = 0.03974517324292037
$ws.Cells.Item(8, 4).Value = 0.1431084544871233
$ws.Cells.Item(8, 5).Value = 0.0609719324907374
$ws.Cells.Item(8, 6).Value = 2.612048684312441
$ws.Cells.Item(8, 8).Value = 0.07973214163530429
$ws.Cells.Item(8, 9).Value = 2.135687664364198
$ws.Cells.Item(8, 11).Value = 0.6056238609952516
$ws.Cells.Item(8, 12).Value = 0.2431844276394912

$ws.Cells.Item(9, 2).Value = 0.839340416755249
$ws.Cells.Item(9, 3).Value = 0.05322382392665759
$ws.Cells.Item(9, 4).Value = 0.1447517514735388
$ws.Cells.Item(9, 5).Value = 0.06223435567550339
$ws.Cells.Item(9, 6).Value = 2.750534342411299
$ws.Cells.Item(9, 8).Value = 0.07973214163530429
$ws.Cells.Item(9, 9).Value = 2.221657233728891
$ws.Cells.Item(9, 11).Value = 0.6858677430425359
$ws.Cells.Item(9, 12).Value = 0.261853296654138

$ws.Cells.Item(10, 2).Value = 0.892488089071179
$ws.Cells.Item(10, 3).Value = 0.06312273449911743
$ws.Cells.Item(10, 4).Value = 0.1459026043700362
$ws.Cells.Item(10, 5).Value = 0.06328740268803301
$ws.Cells.Item(10, 6).Value = 2.856429242130048
$ws.Cells.Item(10, 8).Value = 0.07973214163530429
$ws.Cells.Item(10, 9).Value = 2.287554572073915
$ws.Cells.Item(10, 11).Value = 0.7472251482338947
$ws.Cells.Item(10, 12).Value = 0.2762844874468158

$ws.Cells.Item(11, 2).Value = 0.9173216697803639
$ws.Cells.Item(11, 3).Value = 0.06762696470323704
$ws.Cells.Item(11, 4).Value = 0.1464142668619246
$ws.Cells.Item(11, 5).Value = 0.06379375755386008
$ws.Cells.Item(11, 6).Value = 2.905520194086677
$ws.Cells.Item(11, 8).Value = 0.07973214163530429
$ws.Cells.Item(11, 9).Value = 2.318137281930532
$ws.Cells.Item(11, 11).Value = 0.7756664249608036
$ws.Cells.Item(11, 12).Value = 0.2830065516543101

$ws.Cells.Item(12, 2).Value = 0.9268200931945501
$ws.Cells.Item(12, 3).Value = 0.06933289356206274
$ws.Cells.Item(12, 4).Value = 0.1466063422332979
$ws.Cells.Item(12, 5).Value = 0.06398942885915204
$ws.Cells.Item(12, 6).Value = 2.924242794660472
$ws.Cells.Item(12, 8).Value = 0.07973214163530429
$ws.Cells.Item(12, 9).Value = 2.329805929395675
$ws.Cells.Item(12, 11).Value = 0.7865129249301503
$ws.Cells.Item(12, 12).Value = 0.2855747316202724

$ws.Cells.Item(13, 2).Value = 0.924770233963045
$ws.Cells.Item(13, 3).Value = 0.06896547733187219
$ws.Cells.Item(13, 4).Value = 0.1465650496718425
$ws.Cells.Item(13, 5).Value = 0.06394711298622013
$ws.Cells.Item(13, 6).Value = 2.920204623134794
$ws.Cells.Item(13, 8).Value = 0.07973214163530429
$ws.Cells.Item(13, 9).Value = 2.327288970937957
$ws.Cells.Item(13, 11).Value = 0.7841735345880068
$ws.Cells.Item(13, 12).Value = 0.2850206182951496

$ws.Cells.Item(14, 2).Value = 0.9181012166129108
$ws.Cells.Item(14, 3).Value = 0.06776730637680828
$ws.Cells.Item(14, 4).Value = 0.1464301025402008
$ws.Cells.Item(14, 5).Value = 0.06380977688749923
$ws.Cells.Item(14, 6).Value = 2.907057845161461
$ws.Cells.Item(14, 8).Value = 0.07973214163530429
$ws.Cells.Item(14, 9).Value = 2.319095509045155
$ws.Cells.Item(14, 11).Value = 0.7765572401980592
$ws.Cells.Item(14, 12).Value = 0.2832173824211708

$ws.Cells.Item(15, 2).Value = 0.9140285585467325
$ws.Cells.Item(15, 3).Value = 0.06703343102664405
$ws.Cells.Item(15, 4).Value = 0.146347225503046
$ws.Cells.Item(15, 5).Value = 0.06372616572866363
$ws.Cells.Item(15, 6).Value = 2.899022398192614
$ws.Cells.Item(15, 8).Value = 0.07973214163530429
$ws.Cells.Item(15, 9).Value = 2.31408820639686
$ws.Cells.Item(15, 11).Value = 0.7719019983123019
$ws.Cells.Item(15, 12).Value = 0.2821158040760423

$ws.Cells.Item(16, 2).Value = 0.8908783728921321
$ws.Cells.Item(16, 3).Value = 0.06282840504803744
$ws.Cells.Item(16, 4).Value = 0.14586892948914
$ws.Cells.Item(16, 5).Value = 0.06325486069341935
$ws.Cells.Item(16, 6).Value = 2.853239590081273
$ws.Cells.Item(16, 8).Value = 0.07973214163530429
$ws.Cells.Item(16, 9).Value = 2.285568154236429
$ws.Cells.Item(16, 11).Value = 0.7453771210736306
$ws.Cells.Item(16, 12).Value = 0.275848354715464

$ws.Cells.Item(17, 2).Value = 0.8768446746982193
$ws.Cells.Item(17, 3).Value = 0.06024913228154105
$ws.Cells.Item(17, 4).Value = 0.1455724898620403
$ws.Cells.Item(17, 5).Value = 0.06297272529240416
$ws.Cells.Item(17, 6).Value = 2.825389219166141
$ws.Cells.Item(17, 8).Value = 0.07973214163530429
$ws.Cells.Item(17, 9).Value = 2.268227540585656
$ws.Cells.Item(17, 11).Value = 0.7292407685752096
$ws.Cells.Item(17, 12).Value = 0.272043799660679

$ws.Cells.Item(18, 2).Value = 0.868834643755207
$ws.Cells.Item(18, 3).Value = 0.0587657061134621
$ws.Cells.Item(18, 4).Value = 0.1454008671595055
$ws.Cells.Item(18, 5).Value = 0.06281301996665789
$ws.Cells.Item(18, 6).Value = 2.809456859245699
$ws.Cells.Item(18, 8).Value = 0.07973214163530429
$ws.Cells.Item(18, 9).Value = 2.2583106551225
$ws.Cells.Item(18, 11).Value = 0.7200093865928068
$ws.Cells.Item(18, 12).Value = 0.2698703164732876

$ws.Cells.Item(19, 2).Value = 0.8661331932642895
$ws.Cells.Item(19, 3).Value = 0.05826345800288379
$ws.Cells.Item(19, 4).Value = 0.1453425657617338
$ws.Cells.Item(19, 5).Value = 0.06275938822795268
$ws.Cells.Item(19, 6).Value = 2.804077257378367
$ws.Cells.Item(19, 8).Value = 0.07973214163530429
$ws.Cells.Item(19, 9).Value = 2.254962737428784
$ws.Cells.Item(19, 11).Value = 0.7168923481371507
$ws.Cells.Item(19, 12).Value = 0.2691369516148825

$ws.Cells.Item(20, 2).Value = 0.878332191103965
$ws.Cells.Item(20, 3).Value = 0.06052368880403947
$ws.Cells.Item(20, 4).Value = 0.1456041619338535
$ws.Cells.Item(20, 5).Value = 0.06300249296715421
$ws.Cells.Item(20, 6).Value = 2.828344987096415
$ws.Cells.Item(20, 8).Value = 0.07973214163530429
$ws.Cells.Item(20, 9).Value = 2.270067577721733
$ws.Cells.Item(20, 11).Value = 0.7309533526448888
$ws.Cells.Item(20, 12).Value = 0.2724472691912467

$ws.Cells.Item(21, 2).Value = 0.9200575035957286
$ws.Cells.Item(21, 3).Value = 0.06811922991884956
$ws.Cells.Item(21, 4).Value = 0.1464697851734371
$ws.Cells.Item(21, 5).Value = 0.0638500093065737
$ws.Cells.Item(21, 6).Value = 2.910915758685007
$ws.Cells.Item(21, 8).Value = 0.07973214163530429
$ws.Cells.Item(21, 9).Value = 2.321499744576869
$ws.Cells.Item(21, 11).Value = 0.7787922546839923
$ws.Cells.Item(21, 12).Value = 0.2837464204200018

$ws.Cells.Item(22, 2).Value = 0.9478781889637276
$ws.Cells.Item(22, 3).Value = 0.07308500623241798
$ws.Cells.Item(22, 4).Value = 0.1470257391166498
$ws.Cells.Item(22, 5).Value = 0.06442679114836025
$ws.Cells.Item(22, 6).Value = 2.965655678686005
$ws.Cells.Item(22, 8).Value = 0.07973214163530429
$ws.Cells.Item(22, 9).Value = 2.355624719915269
$ws.Cells.Item(22, 11).Value = 0.8105032238775038
$ws.Cells.Item(22, 12).Value = 0.2912632992868822

$ws.Cells.Item(23, 2).Value = 0.9329793367512877
$ws.Cells.Item(23, 3).Value = 0.07043449129045598
$ws.Cells.Item(23, 4).Value = 0.1467299023511686
$ws.Cells.Item(23, 5).Value = 0.0641168590492569
$ws.Cells.Item(23, 6).Value = 2.93636876267874
$ws.Cells.Item(23, 8).Value = 0.07973214163530429
$ws.Cells.Item(23, 9).Value = 2.337364634249909
$ws.Cells.Item(23, 11).Value = 0.7935376308480215
$ws.Cells.Item(23, 12).Value = 0.2872392773049199

$ws.Cells.Item(24, 2).Value = 0.8776595038806079
$ws.Cells.Item(24, 3).Value = 0.0603995635784571
$ws.Cells.Item(24, 4).Value = 0.1455898467156516
$ws.Cells.Item(24, 5).Value = 0.06298902722027577
$ws.Cells.Item(24, 6).Value = 2.827008437802931
$ws.Cells.Item(24, 8).Value = 0.07973214163530429
$ws.Cells.Item(24, 9).Value = 2.26923553354689
$ws.Cells.Item(24, 11).Value = 0.7301789513032588
$ws.Cells.Item(24, 12).Value = 0.2722648176187619

$ws.Cells.Item(25, 2).Value = 0.8203495690272575
$ws.Cells.Item(25, 3).Value = 0.04957875338688211
$ws.Cells.Item(25, 4).Value = 0.1443172394085011
$ws.Cells.Item(25, 5).Value = 0.06187079730118583
$ws.Cells.Item(25, 6).Value = 2.712346594024012
$ws.Cells.Item(25, 8).Value = 0.07973214163530429
$ws.Cells.Item(25, 9).Value = 2.197923114728653
$ws.Cells.Item(25, 11).Value = 0.6637398292863281
$ws.Cells.Item(25, 12).Value = 0.2566777937795592
